$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.145.04"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.194.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.88%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "80.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.507"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.06%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.465"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0765"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "28.86"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.70"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -11.37%  "
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.540.18"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.17"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.83"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.203.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.706"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -6.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.067.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0865"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.67"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "224.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.92%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.38"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.78"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.39"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.31%  "
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "148.63"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.40"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -10.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.18%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.33"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0689"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.109"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.13"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0948"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.60"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.63"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.57"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.897.35"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.07"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -8.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0258"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.98"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.77"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -10.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.58"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.409.01"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.92%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.95"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.72"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.65%  "
